# Slide 10 ("Interactions") gets a new caption textbox crediting the
# source of the MVT diagram picture already on the slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$shp = $s.Shapes.AddTextbox(1, 470.78741455078125, 359.1378173828125, 205.15740966796875, 26.10236358642578)
$shp.Name = "Shape 89"

# No fill / no outline box.
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

# Body text-box insets + anchoring matching the rest of the deck.
$shp.TextFrame.MarginLeft = 7.198818897637795
$shp.TextFrame.MarginRight = 7.198818897637795
$shp.TextFrame.MarginTop = 7.198818897637795
$shp.TextFrame.MarginBottom = 7.198818897637795
$shp.TextFrame.AutoSize = 0
$shp.TextFrame2.VerticalAnchor = 1
$shp.TextFrame2.HorizontalAnchor = 0

$tr = $shp.TextFrame.TextRange
$tr.Text = "(http://littlegreenriver.com/weblog/wp-content/uploads/mtv-diagram-730x1024.png)"
$tr.ParagraphFormat.SpaceBefore = 0
$tr.ParagraphFormat.Bullet.Visible = 0
$tr.Font.Size = 8
$tr.Font.Color.RGB = 13421772
$tr.LanguageID = "en"
